# Add two new columns (I = "I0", J = "IF") to Sheet1, matching the
# header style already used by the existing header row and filling in
# the per-row numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (font, border, alignment) from the existing header
# cell H1 so the new headers look consistent with the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data ----------------------------------------------------------------
# row, I0, IF
$rows = @(
    @(2, 9, 9),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 8, 8),
    @(6, 7, 7),
    @(7, 6, 6),
    @(8, 9, 9),
    @(9, 7, 7),
    @(10, 7, 7),
    @(11, 8, 8),
    @(12, 8, 8),
    @(13, 8, 8),
    @(14, 10, 10),
    @(15, 6, 6),
    @(16, 7, 7),
    @(17, 8, 8),
    @(18, 9, 9),
    @(19, 7, 7),
    @(20, 7, 7),
    @(21, 6, 6),
    @(22, 7, 7),
    @(23, 6, 6),
    @(24, 6, 6),
    @(25, 7, 7),
    @(26, 6, 6),
    @(27, 8, 8),
    @(28, 6, 6),
    @(29, 6, 6),
    @(30, 5, 6),
    @(31, 7, 7),
    @(32, 7, 7),
    @(33, 6, 6),
    @(34, 6, 7),
    @(35, 8, 8),
    @(36, 6, 6),
    @(37, 6, 6),
    @(38, 7, 7),
    @(39, 7, 7),
    @(40, 7, 7),
    @(41, 7, 7),
    @(42, 6, 6),
    @(43, 5, 5),
    @(44, 6, 6),
    @(45, 6, 6),
    @(46, 5, 5),
    @(47, 3, 3),
    @(48, 4, 4)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 9).Value = $r[1]   # column I
    $ws.Cells.Item($rowNum, 10).Value = $r[2]  # column J
}
